# Ran code for averaged intensities on spiral schemes.
# This reorders the "Gaussian-Quadrature" row to directly follow the "Ring
# Perpendicular to TD" rows, inserts three new "Spiral-*" scheme rows after it,
# and appends the previously-existing NoRotation/Rotation/HexGrid rows below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatted (bordered / bold / centered) row-index column A
# formatting down to the 3 newly added rows (17-19) by copying it from row 16.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10: Gaussian-Quadrature (moved up from old row 16)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.001089611008914
$ws.Range("D10").Value = 0.966172179807626
$ws.Range("E10").Value = 1.002280705661392
$ws.Range("F10").Value = 1.001089611008914
$ws.Range("G10").Value = 0.9837831208736161
$ws.Range("H10").Value = 1.013900055707453
$ws.Range("I10").Value = 1.003413304056899
$ws.Range("J10").Value = 0.966172179807626
$ws.Range("K10").Value = 0.9842264427345089
$ws.Range("L10").Value = 0.9926580268717113
$ws.Range("M10").Value = 0.9951064961859831

# Row 11: Spiral-90deg-10rot-5space (new)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9978672936785771
$ws.Range("D11").Value = 0.9604644855454775
$ws.Range("E11").Value = 1.00665493896581
$ws.Range("F11").Value = 0.9978672936785771
$ws.Range("G11").Value = 0.9757049890430439
$ws.Range("H11").Value = 1.025636454814496
$ws.Range("I11").Value = 1.004927345411139
$ws.Range("J11").Value = 0.9604644855454775
$ws.Range("K11").Value = 0.9835597122556439
$ws.Range("L11").Value = 0.9907135029671104
$ws.Range("M11").Value = 0.9952092512430905

# Row 12: Spiral-90deg-15rot-5space (new)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9977209437614809
$ws.Range("D12").Value = 0.9611979808109042
$ws.Range("E12").Value = 1.006490989386609
$ws.Range("F12").Value = 0.9977209437614809
$ws.Range("G12").Value = 0.9760425390899916
$ws.Range("H12").Value = 1.025113174179572
$ws.Range("I12").Value = 1.004755388229433
$ws.Range("J12").Value = 0.9611979808109042
$ws.Range("K12").Value = 0.9838444850987567
$ws.Range("L12").Value = 0.9907827144301187
$ws.Range("M12").Value = 0.9952201692429985

# Row 13: Spiral-90deg-10rot-3space (new)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9978285136637507
$ws.Range("D13").Value = 0.9606300162015248
$ws.Range("E13").Value = 1.006627111908095
$ws.Range("F13").Value = 0.9978285136637507
$ws.Range("G13").Value = 0.9757770950205678
$ws.Range("H13").Value = 1.025469807010397
$ws.Range("I13").Value = 1.0048919286526
$ws.Range("J13").Value = 0.9606300162015248
$ws.Range("K13").Value = 0.9836285640548099
$ws.Range("L13").Value = 0.9907285388592804
$ws.Range("M13").Value = 0.9952040787428226

# Row 14: NoRotation-tilt60deg (shifted down from old row 10)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.00252
$ws.Range("D14").Value = 0.9563959999999999
$ws.Range("E14").Value = 1.005636000000001
$ws.Range("F14").Value = 1.00252
$ws.Range("G14").Value = 0.9755880000000007
$ws.Range("H14").Value = 1.036508
$ws.Range("I14").Value = 1.004599999999999
$ws.Range("J14").Value = 0.9563959999999999
$ws.Range("K14").Value = 0.9810160000000003
$ws.Range("L14").Value = 0.991768
$ws.Range("M14").Value = 0.9968746666666668

# Row 15: Rotation-NoTilt (shifted down from old row 11)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.02
$ws.Range("D15").Value = 0.8943125
$ws.Range("E15").Value = 1.017662499999999
$ws.Range("F15").Value = 1.02
$ws.Range("G15").Value = 0.9386124999999987
$ws.Range("H15").Value = 1.09
$ws.Range("I15").Value = 1.02
$ws.Range("J15").Value = 0.8943125
$ws.Range("K15").Value = 0.9559874999999995
$ws.Range("L15").Value = 0.9879937499999997
$ws.Range("M15").Value = 0.9967645833333328

# Row 16: Rotation-60detTilt (shifted down from old row 12)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.010254794547201
$ws.Range("D16").Value = 0.9369929713664048
$ws.Range("E16").Value = 1.008972639334405
$ws.Range("F16").Value = 1.010254794547201
$ws.Range("G16").Value = 0.962624454758402
$ws.Range("H16").Value = 1.048998440345594
$ws.Range("I16").Value = 1.010047546163203
$ws.Range("J16").Value = 0.9369929713664048
$ws.Range("K16").Value = 0.9729828053504046
$ws.Range("L16").Value = 0.9916187999488029
$ws.Range("M16").Value = 0.9963151410858683

# Row 17: HexGrid-90degTilt5degRes (shifted down from old row 13)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9961198445830963
$ws.Range("D17").Value = 0.9960265009822649
$ws.Range("E17").Value = 0.9965385110210615
$ws.Range("F17").Value = 0.9961198445830963
$ws.Range("G17").Value = 0.9956908311575564
$ws.Range("H17").Value = 0.9961067649005235
$ws.Range("I17").Value = 0.9962290122209874
$ws.Range("J17").Value = 0.9960265009822649
$ws.Range("K17").Value = 0.9962825060016631
$ws.Range("L17").Value = 0.9962011752923798
$ws.Range("M17").Value = 0.9961185774775817

# Row 18: HexGrid-90degTilt22p5degRes (shifted down from old row 14)
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9966864287182757
$ws.Range("D18").Value = 0.9990765350099139
$ws.Range("E18").Value = 0.9948939286529779
$ws.Range("F18").Value = 0.9966864287182757
$ws.Range("G18").Value = 0.9984206950134585
$ws.Range("H18").Value = 0.9927603304513295
$ws.Range("I18").Value = 0.9949969447992583
$ws.Range("J18").Value = 0.9990765350099139
$ws.Range("K18").Value = 0.9969852318314458
$ws.Range("L18").Value = 0.9968358302748608
$ws.Range("M18").Value = 0.9961391437742023

# Row 19: HexGrid-60degTilt5degRes (shifted down from old row 15)
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9941938876192412
$ws.Range("D19").Value = 1.007739645320272
$ws.Range("E19").Value = 0.9934528830580236
$ws.Range("F19").Value = 0.9941938876192412
$ws.Range("G19").Value = 1.002375766322493
$ws.Range("H19").Value = 0.9878523205576389
$ws.Range("I19").Value = 0.9932433266612383
$ws.Range("J19").Value = 1.007739645320272
$ws.Range("K19").Value = 1.000596264189148
$ws.Range("L19").Value = 0.9973950759041946
$ws.Range("M19").Value = 0.9964763049231512
